# Fix: Updating dashboard. Started Send mail
#
# - Shrink the saved window size recorded for the workbook.
# - Add six new "Mail*" settings (Name/Value pairs) to the Assets sheet,
#   rows 10-15, which also appends their text to the shared-string table.
# - Row 1000 (a trailing, otherwise-empty formatted row) is removed.
# - Leave the active selection on row 11 (the first of the newly added
#   rows), matching where the user ended up after typing the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# Best-effort: record the (now smaller) Excel window size on the workbook.
try { $excel.ActiveWindow.Width = 9750 } catch {}
try { $excel.ActiveWindow.Height = 9795 } catch {}

# New Mail* settings added to the Assets sheet (Name column == Value column).
$mailSettings = @(
    "MailBodyText",
    "MailSenderAddress",
    "MailSenderName",
    "MailServerAddress",
    "MailServerPort",
    "MailSubject"
)

$row = 10
foreach ($name in $mailSettings) {
    $ws.Range("A$row").Value = $name
    $ws.Range("B$row").Value = $name
    $row = $row + 1
}

# Drop the trailing placeholder row 1000 (dimension shrinks by one row).
[void]$ws.Rows("1000:1000").Delete()

# Leave the selection where the user finished editing: the whole of row 11.
[void]$ws.Rows("11:11").Select()
